$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The source data had a bug (convertJsonToExcel) where the
# "last_edited_time" shared string was duplicated across rows 6-13
# (column D). The fix corrects that shared value for every row that
# referenced it.
$ws.Range("D6:D13").Value = "2024-08-24T20:33:00.000Z"

# Update numeric values for row 6 (Tháng 8)
$ws.Range("T6").Value = 104000000
$ws.Range("W6").Value = 187715000
$ws.Range("AA6").Value = 277735000
$ws.Range("AE6").Value = 465450000
$ws.Range("AH6").Value = 377450000
$ws.Range("AK6").Value = 59
$ws.Range("AN6").Value = 88000000
$ws.Range("AQ6").Value = 481450000
